# Applies the Coinranking "cryptos" price/volume refresh captured in the
# commit diff: updated Price (D) / Volume(1h) (E) figures for most rows,
# plus a 4-way reshuffle of the Cosmos/Toncoin/Kaspa/InjectiveProtocol rows
# (28-31) which swapped Coin name, Link, Price and Volume together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => ordered hashtable of column-letter => new text value.
# Using an ordered dictionary keeps B/C/D/E writes in a predictable order.
$rowUpdates = [ordered]@{
    2 = [ordered]@{ "D" = "47.325.29"; "E" = "  +1.05%  " }
    3 = [ordered]@{ "D" = "2.492.89"; "E" = "  +0.44%  " }
    4 = [ordered]@{ "D" = "1.00"; "E" = "  +0.10%  " }
    5 = [ordered]@{ "D" = "320.92"; "E" = "  -0.55%  " }
    6 = [ordered]@{ "D" = "108.63"; "E" = "  +4.28%  " }
    7 = [ordered]@{ "D" = "0.522"; "E" = "  +0.42%  " }
    8 = [ordered]@{ "E" = "  -0.03%  " }
    9 = [ordered]@{ "D" = "0.536"; "E" = "  -0.24%  " }
    10 = [ordered]@{ "D" = "39.15"; "E" = "  +6.35%  " }
    11 = [ordered]@{ "E" = "  -0.07%  " }
    12 = [ordered]@{ "D" = "0.124"; "E" = "  +0.29%  " }
    13 = [ordered]@{ "D" = "18.38"; "E" = "  +1.03%  " }
    14 = [ordered]@{ "E" = "  -0.26%  " }
    15 = [ordered]@{ "D" = "2.882.23"; "E" = "  +0.40%  " }
    16 = [ordered]@{ "D" = "2.497.32"; "E" = "  -1.05%  " }
    17 = [ordered]@{ "D" = "0.846"; "E" = "  +0.86%  " }
    18 = [ordered]@{ "D" = "47.225.39"; "E" = "  +1.01%  " }
    19 = [ordered]@{ "D" = "13.06"; "E" = "  +4.32%  " }
    20 = [ordered]@{ "D" = "6.62"; "E" = "  +0.85%  " }
    21 = [ordered]@{ "D" = "0.0₃0935"; "E" = "  +0.64%  " }
    22 = [ordered]@{ "D" = "2.65"; "E" = "  +13.15%  " }
    23 = [ordered]@{ "E" = "  -0.29%  " }
    24 = [ordered]@{ "D" = "245.17"; "E" = "  -1.90%  " }
    25 = [ordered]@{ "D" = "2.55"; "E" = "  +0.69%  " }
    26 = [ordered]@{ "E" = "  +0.21%  " }
    27 = [ordered]@{ "E" = "  -1.02%  " }
    28 = [ordered]@{ "B" = "Cosmos"; "C" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; "D" = "10.01"; "E" = "  -0.13%  " }
    29 = [ordered]@{ "B" = "Toncoin"; "C" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; "D" = "2.19"; "E" = "  +0.04%  " }
    30 = [ordered]@{ "B" = "Kaspa"; "C" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; "D" = "0.137"; "E" = "  +3.96%  " }
    31 = [ordered]@{ "B" = "InjectiveProtocol"; "C" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; "D" = "34.76"; "E" = "  -0.64%  " }
    32 = [ordered]@{ "D" = "49.71"; "E" = "  +0.65%  " }
    33 = [ordered]@{ "D" = "20.62"; "E" = "  +5.68%  " }
    34 = [ordered]@{ "D" = "5.36"; "E" = "  +1.18%  " }
    35 = [ordered]@{ "D" = "0.0784"; "E" = "  +1.35%  " }
    36 = [ordered]@{ "E" = "  +0.09%  " }
    37 = [ordered]@{ "D" = "4.73"; "E" = "  +4.10%  " }
    38 = [ordered]@{ "E" = "  +2.92%  " }
    39 = [ordered]@{ "D" = "2.94"; "E" = "  -0.56%  " }
    40 = [ordered]@{ "D" = "22.91"; "E" = "  +7.53%  " }
    41 = [ordered]@{ "E" = "  -0.02%  " }
    42 = [ordered]@{ "E" = "  +0.22%  " }
    43 = [ordered]@{ "D" = "116.03"; "E" = "  -5.31%  " }
    44 = [ordered]@{ "E" = "  +0.95%  " }
    45 = [ordered]@{ "D" = "1.996.82"; "E" = "  +2.54%  " }
    46 = [ordered]@{ "D" = "3.03"; "E" = "  +2.47%  " }
    47 = [ordered]@{ "D" = "2.00"; "E" = "  -5.39%  " }
    48 = [ordered]@{ "D" = "9.14"; "E" = "  +0.31%  " }
    49 = [ordered]@{ "E" = "  -0.60%  " }
    50 = [ordered]@{ "D" = "5.10"; "E" = "  -4.67%  " }
    51 = [ordered]@{ "D" = "56.65"; "E" = "  +4.35%  " }
}

foreach ($row in $rowUpdates.Keys) {
    foreach ($col in $rowUpdates[$row].Keys) {
        $newValue = $rowUpdates[$row][$col]
        $cell = $ws.Range("$col$row")
        if ($col -eq "D") {
            # Price column holds plain text (e.g. "47.325.29", "1.00", "0.137");
            # force text storage so Excel does not reinterpret/round it as a number,
            # then clear the temporary number format so no style index is left behind.
            $cell.NumberFormat = "@"
            $cell.Value = $newValue
            $cell.ClearFormats()
        } else {
            $cell.Value = $newValue
        }
    }
}
